$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.060.14"
$ws.Range("D3").Value = "2.542.53"
$ws.Range("E3").Value = "  -0.04%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "318.27"
$ws.Range("E5").Value = "  +0.38%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "97.08"
$ws.Range("E6").Value = "  +2.30%  "
$ws.Range("E7").Value = "  -0.25%  "
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.541"
$ws.Range("E9").Value = "  +1.29%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.22"
$ws.Range("E10").Value = "  -0.41%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0821"
$ws.Range("E11").Value = "  +0.72%  "
$ws.Range("E12").Value = "  +0.14%  "
$ws.Range("E13").Value = "  -4.32%  "
$ws.Range("D14").Value = "2.932.36"
$ws.Range("E14").Value = "  -0.04%  "
$ws.Range("D15").Value = "2.508.64"
$ws.Range("E15").Value = "  -1.68%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "15.16"
$ws.Range("E16").Value = "  -3.56%  "
$ws.Range("E17").Value = "  -1.69%  "
$ws.Range("D18").Value = "43.082.69"
$ws.Range("E18").Value = "  +0.97%  "
$ws.Range("E19").Value = "  +3.29%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.67"
$ws.Range("E20").Value = "  -2.98%  "
$ws.Range("D21").Value = "0.0₃0969"
$ws.Range("E21").Value = "  +0.13%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "69.89"
$ws.Range("E22").Value = "  -1.58%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "254.68"
$ws.Range("E23").Value = "  +0.72%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.98"
$ws.Range("E24").Value = "  +0.31%  "
$ws.Range("E25").Value = "  +2.24%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "26.59"
$ws.Range("E26").Value = "  -2.90%  "
$ws.Range("E27").Value = "  -0.06%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "41.37"
$ws.Range("E28").Value = "  +5.74%  "
$ws.Range("E29").Value = "  +1.97%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "10.41"
$ws.Range("E30").Value = "  +1.99%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "5.92"
$ws.Range("E31").Value = "  -0.73%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "156.51"
$ws.Range("E32").Value = "  +0.67%  "
$ws.Range("E33").Value = "  -1.21%  "
$ws.Range("B34").Value = "LidoDAOToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.36"
$ws.Range("E34").Value = "  -0.43%  "
$ws.Range("B35").Value = "Celestia"
$ws.Range("C35").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "19.20"
$ws.Range("E35").Value = "  -0.75%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.70"
$ws.Range("E36").Value = "  +3.31%  "
$ws.Range("E37").Value = "  +2.06%  "
$ws.Range("E38").Value = "  +1.36%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.46"
$ws.Range("E39").Value = "  +5.47%  "
$ws.Range("E40").Value = "  -0.26%  "
$ws.Range("E41").Value = "  -8.75%  "
$ws.Range("E42").Value = "  -0.24%  "
$ws.Range("E43").Value = "  +1.16%  "
$ws.Range("E44").Value = "  -0.01%  "
$ws.Range("E45").Value = "  -1.70%  "
$ws.Range("D46").Value = "2.003.45"
$ws.Range("E46").Value = "  -1.88%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.13"
$ws.Range("E47").Value = "  +2.06%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "84.79"
$ws.Range("E48").Value = "  +0.23%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "106.15"
$ws.Range("E49").Value = "  +4.10%  "
$ws.Range("B50").Value = "ordi"
$ws.Range("C50").Value = "https://coinranking.com/coin/j7-7vPrOi+ordi-ordi"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "74.99"
$ws.Range("E50").Value = "  +1.03%  "
$ws.Range("B51").Value = "RocketPoolETH"
$ws.Range("C51").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D51").Value = "2.786.38"
$ws.Range("E51").Value = "  -0.04%  "
